# A new weekly price record was inserted for row 414 (pushing the
# existing rows 414-471 down to 415-472). Use a native row insert so
# every subsequent row (and its formatting) shifts down automatically,
# then populate the newly blank row 414 with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 414; rows 414..471 shift to 415..472.
$ws.Rows.Item(414).Insert()

# Fill in the new row 414 with the new weekly record.
$ws.Range("A414").Value = 5
$ws.Range("B414").Value = "Macroferia Regional de Talca"
$ws.Range("C414").Value = "Maule"
$ws.Range("D414").Value = 45131
$ws.Range("E414").Value = 7
$ws.Range("F414").Value = 100112009
$ws.Range("G414").Value = "Acelga"
$ws.Range("H414").Value = "Sin especificar"
$ws.Range("I414").Value = "Primera"
$ws.Range("J414").Value = 500
$ws.Range("K414").Value = 1600
$ws.Range("L414").Value = 1600
$ws.Range("M414").Value = 1600
$ws.Range("N414").Value = "$/docena de atados (4 kilos)"
$ws.Range("O414").Value = "Región del Maule"
$ws.Range("P414").Value = 400
$ws.Range("Q414").Value = 4
$ws.Range("R414").Value = "Hortaliza"
